# Adds line_items.tsv field documentation rows (32-49) to the data dictionary,
# matching the new columns added to line_items.tsv (incl. metric_id join key).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "line_items.tsv"
$ws.Range("B32").Value = "metric_id"
$ws.Range("C32").Value = "A unique ID associated with the specified metric, which can be used to join to the metrics.tsv table"
$ws.Rows.Item(32).RowHeight = 90

$ws.Range("A33").Value = "line_items.tsv"
$ws.Range("B33").Value = "indicator"
$ws.Range("C33").Value = "The indicator associated with the metric"
$ws.Rows.Item(33).RowHeight = 90

$ws.Range("A34").Value = "line_items.tsv"
$ws.Range("B34").Value = "score"
$ws.Range("C34").Value = "A numeric score that assesses country performance against the metric"
$ws.Rows.Item(34).RowHeight = 90

$ws.Range("A35").Value = "line_items.tsv"
$ws.Range("B35").Value = "attribute"
$ws.Range("C35").Value = "The attribute that is required to obtain the specified score on the metric"
$ws.Rows.Item(35).RowHeight = 90

$ws.Range("A36").Value = "line_items.tsv"
$ws.Range("B36").Value = "requirement"
$ws.Range("C36").Value = "The requirement specified based on the attribute"
$ws.Rows.Item(36).RowHeight = 90

$ws.Range("A37").Value = "line_items.tsv"
$ws.Range("B37").Value = "activity"
$ws.Range("C37").Value = "The activity required to address the requirement"
$ws.Rows.Item(37).RowHeight = 90

$ws.Range("A38").Value = "line_items.tsv"
$ws.Range("B38").Value = "unit_cost"
$ws.Range("C38").Value = "The unit cost associated with the activity"
$ws.Rows.Item(38).RowHeight = 90

$ws.Range("A39").Value = "line_items.tsv"
$ws.Range("B39").Value = "unit"
$ws.Range("C39").Value = "The units associated with that unit cost"
$ws.Rows.Item(39).RowHeight = 90

$ws.Range("A40").Value = "line_items.tsv"
$ws.Range("B40").Value = "description"
$ws.Range("C40").Value = "A brief written description of the activity"
$ws.Rows.Item(40).RowHeight = 90

$ws.Range("A41").Value = "line_items.tsv"
$ws.Range("B41").Value = "administrative_level"
$ws.Range("C41").Value = "Defines country geopolitical organization. For the purposes of costing IHR implementation, these are the administrative units that support public health efforts such as biosurveillance or emergency response. Administrative organization is divided between intermediate and local levels, including an optional second intermediate level for countries where such an organization exists."
$ws.Range("D41").Value = "Country: central or national-level government`nIntermediate (e.g., province, district): Largest geopolitical unit under the central government`nLocal (e.g., county, city): Smallest geopolitical unit with a role in national public health prevention, detection, and response efforts`nHealth facility: Primarily expected to be hospitals and government-run health centers for the purposes of costing IHR implementation. Includes facilities participating in IHR-related activities including biosurveillance programs, point-of-care diagnostics for priorities diseases, prevention of healthcare associated infections, and biosafety and biosecurity programs.`nPopulation: Population, used to calculate costs for activities that scale with total population size (e.g., cost per vaccine dose)"
$ws.Range("D41").Characters(1, 9).Font.Bold = $true
$ws.Range("D41").Characters(47, 41).Font.Bold = $true
$ws.Range("D41").Characters(143, 28).Font.Bold = $true
$ws.Range("D41").Characters(280, 17).Font.Bold = $true
$ws.Range("D41").Characters(649, 11).Font.Bold = $true
$ws.Rows.Item(41).RowHeight = 106

$ws.Range("A42").Value = "line_items.tsv"
$ws.Range("B42").Value = "cost_type"
$ws.Range("C42").Value = "Indicates whther the cost is a one-time/start up cost or a recurring cost. All recurring costs are assumed to be annual unless another (temporal) custom multiplier is specified."
$ws.Range("D42").Value = "One-time costs are needed once, assumed to be during year 1 or at startup`nRecurring costs are needed multiple times, assumed to be annually unless another (temporal) custom multiplier is specified"
$ws.Range("D42").Characters(1, 9).Font.Bold = $true
$ws.Range("D42").Characters(75, 10).Font.Bold = $true
$ws.Rows.Item(42).RowHeight = 90

$ws.Range("A43").Value = "line_items.tsv"
$ws.Range("B43").Value = "custom_multiplier_1"
$ws.Range("C43").Value = "A custom multiplier that can be used to adjust the cost of the line item, for example, specifying how many times an activity is needed per year, or how many days a meeting is, or how many software licenses are needed."
$ws.Rows.Item(43).RowHeight = 90

$ws.Range("A44").Value = "line_items.tsv"
$ws.Range("B44").Value = "custom_multiplier_1_units"
$ws.Range("C44").Value = "The units associated with the custom multiplier"
$ws.Rows.Item(44).RowHeight = 90

$ws.Range("A45").Value = "line_items.tsv"
$ws.Range("B45").Value = "custom_multiplier_2"
$ws.Range("C45").Value = "A custom multiplier that can be used to adjust the cost of the line item, for example, specifying how many times an activity is needed per year, or how many days a meeting is, or how many software licenses are needed."
$ws.Rows.Item(45).RowHeight = 90

$ws.Range("A46").Value = "line_items.tsv"
$ws.Range("B46").Value = "custom_multiplier_2_units"
$ws.Range("C46").Value = "The units associated with the custom multiplier"
$ws.Rows.Item(46).RowHeight = 90

$ws.Range("A47").Value = "line_items.tsv"
$ws.Range("B47").Value = "relevant_references"
$ws.Range("C47").Value = "Any additional references associated with the designation of the line item"
$ws.Rows.Item(47).RowHeight = 90

$ws.Range("A48").Value = "line_items.tsv"
$ws.Range("B48").Value = "optional_cost"
$ws.Range("C48").Value = "Boolean (TRUE/FALSE) that indicates whether or not the cost is considered an `"optional`" cost based on an interpretation of the language in the JEE or other relevant metric"
$ws.Rows.Item(48).RowHeight = 90

$ws.Range("A49").Value = "line_items.tsv"
$ws.Range("B49").Value = "notes_assumptions"
$ws.Range("C49").Value = "Any notes or additional assumptions made associated with the line item"
$ws.Rows.Item(49).RowHeight = 90

$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 2
$ws.Range("C48").Select()

